$d = $word.ActiveDocument

# --- Change 1: " de prés ilz font" -> " de pres ilz font" -------------
# (the accented "e" becomes a plain "e"; Word naturally splits the run
# around the replaced character when only part of a run's text changes)
$found1 = $d.Content
$found1.Find.Execute(" de prés ilz font", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $found1.Start
$accent1 = $d.Range($start1 + 6, $start1 + 7)
$accent1.Text = "e"

# --- Change 2: "près." -> "pres." ---------------------------------------
$found2 = $d.Content
$found2.Find.Execute("près.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $found2.Start
$accent2 = $d.Range($start2 + 2, $start2 + 3)
$accent2.Text = "e"
